$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Header row: new column H ("Request ID"), formatted like the rest of row 1 ---
$ws.Range("G1").Copy()
$ws.Range("H1").PasteSpecial(-4122)  # xlPasteFormats - reuse the existing header style
$ws.Range("H1").Value = "Request ID"

# --- Row 2 updates ---
# D2 / E2 switch from text to real numbers
$ws.Range("D2").Value = 12
$ws.Range("E2").Value = 20
# F2 / G2 are left untouched (unchanged in the diff)
# H2 is a (practically) blank text cell
$ws.Range("H2").Value = " "

# --- Row 3 (new) ---
$ws.Range("A3").Value = "Transmittance"
$ws.Range("B3").Value = "Tool B"
$ws.Range("C3").Value = "Yes"
$ws.Range("D3").Value = 10
$ws.Range("E3").Value = " "
$ws.Range("F3").Value = "'2025"
$ws.Range("G3").Value = "nopthing"
$ws.Range("H3").Value = 2

# --- Row 4 (new) ---
$ws.Range("A4").Value = "Transmittance"
$ws.Range("B4").Value = "Tool A"
$ws.Range("C4").Value = "Yes"
$ws.Range("D4").Value = "'20"
$ws.Range("E4").Value = "'30"
$ws.Range("F4").Value = "'25"
$ws.Range("G4").Value = "notes"
$ws.Range("H4").Value = 3
